$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C14").Value = 62.5
